$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 656.9231
$ws.Range("I92").Value = 612.7273
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 612.7273
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 635.2727
$ws.Range("N92").Value = -3396
$ws.Range("H113").Value = 2999
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = 255
$ws.Range("H116").Value = 10979.129
$ws.Range("I116").Value = 12293.048
$ws.Range("J116").Value = 8219.9
$ws.Range("K116").Value = 12293.048
$ws.Range("L116").Value = 8219.9
$ws.Range("M116").Value = -8851.048000000001
$ws.Range("N116").Value = -15103.9
$ws.Range("H127").Value = 3091.875
$ws.Range("I127").Value = 3400.7144
$ws.Range("K127").Value = 10202.1432
$ws.Range("M127").Value = -5242.143199999999
$ws.Range("H132").Value = 75407.09
$ws.Range("I132").Value = 79894.07000000001
$ws.Range("K132").Value = 239682.21
$ws.Range("M132").Value = -237152.21
$ws.Range("H138").Value = 2596.2964
$ws.Range("I138").Value = 2055.45
$ws.Range("J138").Value = 4141.5713
$ws.Range("K138").Value = 6166.349999999999
$ws.Range("L138").Value = 12424.7139
$ws.Range("M138").Value = -1026.349999999999
$ws.Range("N138").Value = -22704.7139
$ws.Range("H141").Value = 1334
$ws.Range("I141").Value = 1334
$ws.Range("K141").Value = 4002
$ws.Range("M141").Value = 1178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 152971
$ws.Range("I6").Value = 8599
$ws.Range("K6").Value = 8599
$ws.Range("M6").Value = -8426
$ws.Range("H32").Value = 5267287
$ws.Range("I32").Value = 6252629
$ws.Range("K32").Value = 6252629
$ws.Range("M32").Value = -6252342
$ws.Range("H45").Value = 3917.36
$ws.Range("I45").Value = 3772.1765
$ws.Range("J45").Value = 4225.875
$ws.Range("K45").Value = 3772.1765
$ws.Range("L45").Value = 4225.875
$ws.Range("M45").Value = -3395.1765
$ws.Range("N45").Value = -4979.875
$ws.Range("H61").Value = 1398299.5
$ws.Range("I61").Value = 1597425.2
$ws.Range("K61").Value = 1597425.2
$ws.Range("M61").Value = -1597213.2
$ws.Range("H132").Value = 934310.8
$ws.Range("I132").Value = 1110046.4
$ws.Range("K132").Value = 3330139.2
$ws.Range("M132").Value = -3327609.2
$ws.Range("H136").Value = 1398299.5
$ws.Range("I136").Value = 1597425.2
$ws.Range("K136").Value = 4792275.6
$ws.Range("M136").Value = -4789725.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1793.3334
$ws.Range("I20").Value = 1937.8125
$ws.Range("J20").Value = 1504.375
$ws.Range("K20").Value = 1937.8125
$ws.Range("L20").Value = 1504.375
$ws.Range("M20").Value = -1690.8125
$ws.Range("N20").Value = -1998.375
$ws.Range("H86").Value = 1827.1154
$ws.Range("I86").Value = 1747.5652
$ws.Range("K86").Value = 1747.5652
$ws.Range("M86").Value = -624.5652
$ws.Range("H89").Value = 1827.1154
$ws.Range("I89").Value = 1747.5652
$ws.Range("K89").Value = 8737.826000000001
$ws.Range("M89").Value = -3121.826000000001
$ws.Range("H105").Value = 2199.7334
$ws.Range("I105").Value = 2155.0833
$ws.Range("J105").Value = 2378.3333
$ws.Range("K105").Value = 2155.0833
$ws.Range("L105").Value = 2378.3333
$ws.Range("M105").Value = -408.0832999999998
$ws.Range("N105").Value = -5872.3333
$ws.Range("H134").Value = 1556308.2
$ws.Range("I134").Value = 2069595.4
$ws.Range("K134").Value = 6208786.199999999
$ws.Range("M134").Value = -6206251.199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 185.53334
$ws.Range("J7").Value = 360.14285
$ws.Range("L7").Value = 360.14285
$ws.Range("N7").Value = -586.14285
$ws.Range("H12").Value = 2612.5
$ws.Range("I12").Value = 1816.6666
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 1816.6666
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -1646.6666
$ws.Range("N12").Value = -5340
$ws.Range("H31").Value = 6145
$ws.Range("I31").Value = 1652.4615
$ws.Range("J31").Value = 8684.261
$ws.Range("K31").Value = 1652.4615
$ws.Range("L31").Value = 8684.261
$ws.Range("M31").Value = -1357.4615
$ws.Range("N31").Value = -9274.261
$ws.Range("H34").Value = 6145
$ws.Range("I34").Value = 1652.4615
$ws.Range("J34").Value = 8684.261
$ws.Range("K34").Value = 1652.4615
$ws.Range("L34").Value = 8684.261
$ws.Range("M34").Value = -1450.4615
$ws.Range("N34").Value = -9088.261
$ws.Range("H68").Value = 81976.09
$ws.Range("J68").Value = 81976.09
$ws.Range("L68").Value = 81976.09
$ws.Range("N68").Value = -83474.09
$ws.Range("H71").Value = 81976.09
$ws.Range("J71").Value = 81976.09
$ws.Range("L71").Value = 245928.27
$ws.Range("N71").Value = -253416.27
$ws.Range("H122").Value = 3084.6191
$ws.Range("I122").Value = 1217.9231
$ws.Range("J122").Value = 6118
$ws.Range("K122").Value = 3653.7693
$ws.Range("L122").Value = 18354
$ws.Range("M122").Value = -1203.7693
$ws.Range("N122").Value = -23254

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143466830
$ws.Range("I4").Value = 177166860
$ws.Range("K4").Value = 531500580
$ws.Range("M4").Value = -531500468
$ws.Range("H114").Value = 1349.3103
$ws.Range("I114").Value = 204.38889
$ws.Range("K114").Value = 613.1666700000001
$ws.Range("M114").Value = 2640.83333
$ws.Range("H136").Value = 1565.909
$ws.Range("I136").Value = 1565.909
$ws.Range("K136").Value = 4697.727000000001
$ws.Range("M136").Value = 402.2729999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2591.3914
$ws.Range("I102").Value = 1741.3529
$ws.Range("K102").Value = 1741.3529
$ws.Range("M102").Value = -119.3529000000001
$ws.Range("H132").Value = 804898.25
$ws.Range("I132").Value = 928229.3
$ws.Range("K132").Value = 2784687.9
$ws.Range("M132").Value = -2782157.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2999.875
$ws.Range("I68").Value = 2571.2856
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 2571.2856
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -1822.2856
$ws.Range("N68").Value = -7498
$ws.Range("H71").Value = 2999.875
$ws.Range("I71").Value = 2571.2856
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 12856.428
$ws.Range("L71").Value = 30000
$ws.Range("M71").Value = -9112.428
$ws.Range("N71").Value = -37488
$ws.Range("H100").Value = 7285.048
$ws.Range("I100").Value = 2650.1333
$ws.Range("J100").Value = 18872.334
$ws.Range("K100").Value = 2650.1333
$ws.Range("L100").Value = 18872.334
$ws.Range("M100").Value = -2109.1333
$ws.Range("N100").Value = -19954.334
$ws.Range("H132").Value = 807410
$ws.Range("I132").Value = 1117980.1
$ws.Range("K132").Value = 3353940.3
$ws.Range("M132").Value = -3351410.3
$ws.Range("H136").Value = 3931.9092
$ws.Range("J136").Value = 6251
$ws.Range("L136").Value = 18753
$ws.Range("N136").Value = -23853

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("H113").Value = 3922.9048
$ws.Range("J113").Value = 5999.6
$ws.Range("L113").Value = 17998.8
$ws.Range("N113").Value = -22338.8
$ws.Range("H116").Value = 156761.67
$ws.Range("J116").Value = 156761.67
$ws.Range("L116").Value = 156761.67
$ws.Range("N116").Value = -165939.67
$ws.Range("H132").Value = 5299047
$ws.Range("I132").Value = 6494131.5
$ws.Range("K132").Value = 19482394.5
$ws.Range("M132").Value = -19479864.5
